# Generate Report for Handoff
# Updates the Priority column ("low" -> "ht") and refreshes the
# "Latest Handoff Datetime" timestamp for the files that were just
# handed off, on both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 (the four files whose Priority is still "low")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-31 18:36:46"

# de-de: rows 4-7 (same four files)
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-31 18:36:55"
